$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 770.619
$ws.Range("I12").Value = 111.5625
$ws.Range("K12").Value = 111.5625
$ws.Range("M12").Value = 58.4375
$ws.Range("H32").Value = 10801.7
$ws.Range("I32").Value = 12005.25
$ws.Range("J32").Value = 9999.333000000001
$ws.Range("K32").Value = 12005.25
$ws.Range("L32").Value = 9999.333000000001
$ws.Range("M32").Value = -11679.25
$ws.Range("N32").Value = -10651.333
$ws.Range("H113").Value = 6603.923
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 7168.273
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 7168.273
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -13676.273
$ws.Range("H116").Value = 10575.625
$ws.Range("I116").Value = 10699.75
$ws.Range("K116").Value = 10699.75
$ws.Range("M116").Value = -7257.75
$ws.Range("H121").Value = 2347
$ws.Range("J121").Value = 2347
$ws.Range("L121").Value = 7041
$ws.Range("N121").Value = -10535
$ws.Range("H132").Value = 1106.725
$ws.Range("I132").Value = 916.02856
$ws.Range("J132").Value = 2441.6
$ws.Range("K132").Value = 2748.08568
$ws.Range("L132").Value = 7324.799999999999
$ws.Range("M132").Value = -218.0856800000001
$ws.Range("N132").Value = -12384.8
$ws.Range("H135").Value = 1227.5294
$ws.Range("I135").Value = 991.13336
$ws.Range("J135").Value = 3000.5
$ws.Range("K135").Value = 8920.20024
$ws.Range("L135").Value = 27004.5
$ws.Range("M135").Value = -6385.20024
$ws.Range("N135").Value = -32074.5
$ws.Range("H138").Value = 3823.9773
$ws.Range("I138").Value = 3663.6365
$ws.Range("J138").Value = 3877.4243
$ws.Range("K138").Value = 10990.9095
$ws.Range("L138").Value = 11632.2729
$ws.Range("M138").Value = -5850.9095
$ws.Range("N138").Value = -21912.2729

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8440.808000000001
$ws.Range("I61").Value = 7519.3887
$ws.Range("K61").Value = 7519.3887
$ws.Range("M61").Value = -7307.3887
$ws.Range("H74").Value = 15874881
$ws.Range("I74").Value = 22223840
$ws.Range("K74").Value = 22223840
$ws.Range("M74").Value = -22222966
$ws.Range("H77").Value = 15874881
$ws.Range("I77").Value = 22223840
$ws.Range("K77").Value = 111119200
$ws.Range("M77").Value = -111114832
$ws.Range("H122").Value = 2372.6365
$ws.Range("I122").Value = 1298.3334
$ws.Range("K122").Value = 3895.0002
$ws.Range("M122").Value = -1445.0002
$ws.Range("H136").Value = 8440.808000000001
$ws.Range("I136").Value = 7519.3887
$ws.Range("K136").Value = 22558.1661
$ws.Range("M136").Value = -20008.1661
$ws.Range("H141").Value = 99999
$ws.Range("J141").Value = 99999
$ws.Range("L141").Value = 99999
$ws.Range("N141").Value = -110359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2186.8333
$ws.Range("I99").Value = 2251.75
$ws.Range("J99").Value = 2057
$ws.Range("K99").Value = 2251.75
$ws.Range("L99").Value = 2057
$ws.Range("M99").Value = -753.75
$ws.Range("N99").Value = -5053
$ws.Range("H122").Value = 6028.8335
$ws.Range("I122").Value = 1967.0769
$ws.Range("K122").Value = 5901.2307
$ws.Range("M122").Value = -3451.2307
$ws.Range("H126").Value = 2186.8333
$ws.Range("I126").Value = 2251.75
$ws.Range("J126").Value = 2057
$ws.Range("K126").Value = 6755.25
$ws.Range("L126").Value = 6171
$ws.Range("M126").Value = -4285.25
$ws.Range("N126").Value = -11111
$ws.Range("H141").Value = 172490.62
$ws.Range("J141").Value = 172490.62
$ws.Range("L141").Value = 172490.62
$ws.Range("N141").Value = -182850.62

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4046.0908
$ws.Range("I113").Value = 2870.3333
$ws.Range("J113").Value = 5457
$ws.Range("K113").Value = 2870.3333
$ws.Range("L113").Value = 5457
$ws.Range("M113").Value = -700.3332999999998
$ws.Range("N113").Value = -9797
$ws.Range("H122").Value = 3826.6333
$ws.Range("I122").Value = 3420.1365
$ws.Range("J122").Value = 4944.5
$ws.Range("K122").Value = 10260.4095
$ws.Range("L122").Value = 14833.5
$ws.Range("M122").Value = -7810.4095
$ws.Range("N122").Value = -19733.5
$ws.Range("H132").Value = 5332.793
$ws.Range("I132").Value = 3126.75
$ws.Range("J132").Value = 10235.111
$ws.Range("K132").Value = 9380.25
$ws.Range("L132").Value = 30705.333
$ws.Range("M132").Value = -6850.25
$ws.Range("N132").Value = -35765.333
$ws.Range("H140").Value = 69332.39999999999
$ws.Range("J140").Value = 69332.39999999999
$ws.Range("L140").Value = 69332.39999999999
$ws.Range("N140").Value = -79692.39999999999
$ws.Range("H141").Value = 70429
$ws.Range("J141").Value = 70429
$ws.Range("L141").Value = 70429
$ws.Range("N141").Value = -80789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2314
$ws.Range("I16").Value = 2143.5334
$ws.Range("J16").Value = 3166.3333
$ws.Range("K16").Value = 2143.5334
$ws.Range("L16").Value = 3166.3333
$ws.Range("M16").Value = -1973.5334
$ws.Range("N16").Value = -3506.3333
$ws.Range("H46").Value = 6143
$ws.Range("I46").Value = 3999
$ws.Range("K46").Value = 3999
$ws.Range("M46").Value = -3811
$ws.Range("H55").Value = 2727.3635
$ws.Range("I55").Value = 812.9231
$ws.Range("J55").Value = 5492.6665
$ws.Range("K55").Value = 812.9231
$ws.Range("L55").Value = 5492.6665
$ws.Range("M55").Value = -639.9231
$ws.Range("N55").Value = -5838.6665
$ws.Range("H93").Value = 3982.0833
$ws.Range("I93").Value = 3624
$ws.Range("J93").Value = 4483.4
$ws.Range("K93").Value = 3624
$ws.Range("L93").Value = 4483.4
$ws.Range("M93").Value = -2376
$ws.Range("N93").Value = -6979.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 50000
$ws.Range("I86").Value = 50000
$ws.Range("K86").Value = 50000
$ws.Range("M86").Value = -48877
$ws.Range("H89").Value = 50000
$ws.Range("I89").Value = 50000
$ws.Range("K89").Value = 250000
$ws.Range("M89").Value = -244384
$ws.Range("H132").Value = 6216.4346
$ws.Range("I132").Value = 6645
$ws.Range("J132").Value = 5412.875
$ws.Range("K132").Value = 19935
$ws.Range("L132").Value = 16238.625
$ws.Range("M132").Value = -17405
$ws.Range("N132").Value = -21298.625
